$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.73813591030776
$ws.Range("D2").Value = 4.826604299473929
$ws.Range("E2").Value = 21.47658616257852
$ws.Range("F2").Value = 22.8697344708558
$ws.Range("G2").Value = 26.29328704279751
$ws.Range("H2").Value = 13.73259625338684
$ws.Range("K2").Value = 12.03577553756201
$ws.Range("L2").Value = 7.934078062814927
$ws.Range("M2").Value = 14.58959126765185
$ws.Range("O2").Value = 20.56275605144707

$ws.Range("B3").Value = 15.65424054670321
$ws.Range("D3").Value = 4.749414101993303
$ws.Range("E3").Value = 21.56438718243541
$ws.Range("F3").Value = 22.90814565811468
$ws.Range("G3").Value = 26.35819782191314
$ws.Range("H3").Value = 13.777984373216
$ws.Range("K3").Value = 11.80823693697445
$ws.Range("L3").Value = 7.908617977980042
$ws.Range("M3").Value = 14.56632286585677
$ws.Range("O3").Value = 20.63533945619076

$ws.Range("B4").Value = 15.60535374633824
$ws.Range("D4").Value = 4.700738543284764
$ws.Range("E4").Value = 21.62123997517299
$ws.Range("F4").Value = 22.93786783794481
$ws.Range("G4").Value = 26.40681311736239
$ws.Range("H4").Value = 13.80799037409581
$ws.Range("K4").Value = 11.66514542806026
$ws.Range("L4").Value = 7.893892611023823
$ws.Range("M4").Value = 14.5539385656488
$ws.Range("O4").Value = 20.68430101451257

$ws.Range("B5").Value = 15.58610861314475
$ws.Range("D5").Value = 4.68059359669985
$ws.Range("E5").Value = 21.64514988243929
$ws.Range("F5").Value = 22.9515203666901
$ws.Range("G5").Value = 26.42881835619118
$ws.Range("H5").Value = 13.82075568717442
$ws.Range("K5").Value = 11.60603310876901
$ws.Range("L5").Value = 7.888123868014851
$ws.Range("M5").Value = 14.54937417150281
$ws.Range("O5").Value = 20.70535661215447

$ws.Range("B6").Value = 15.58295430412951
$ws.Range("D6").Value = 4.677230246134451
$ws.Range("E6").Value = 21.64916497383254
$ws.Range("F6").Value = 22.9538803060425
$ws.Range("G6").Value = 26.43260455788187
$ws.Range("H6").Value = 13.82290783273094
$ws.Range("K6").Value = 11.59617062735234
$ws.Range("L6").Value = 7.887180090791188
$ws.Range("M6").Value = 14.54864549576729
$ws.Range("O6").Value = 20.70891946463458

$ws.Range("B7").Value = 15.60509143897885
$ws.Range("D7").Value = 4.700468096529299
$ws.Range("E7").Value = 21.62155942583331
$ws.Range("F7").Value = 22.93804572750809
$ws.Range("G7").Value = 26.40710101593094
$ws.Range("H7").Value = 13.80816035458046
$ws.Range("K7").Value = 11.66435139700923
$ws.Range("L7").Value = 7.89381386760797
$ws.Range("M7").Value = 14.55387505085184
$ws.Range("O7").Value = 20.68458051261921

$ws.Range("B8").Value = 15.70867407136204
$ws.Range("D8").Value = 4.800262358844319
$ws.Range("E8").Value = 21.50625041432967
$ws.Range("F8").Value = 22.88170327241054
$ws.Range("G8").Value = 26.31384528572043
$ws.Range("H8").Value = 13.74780251555916
$ws.Range("K8").Value = 11.95805087765166
$ws.Range("L8").Value = 7.925113486871934
$ws.Range("M8").Value = 14.58117584874903
$ws.Range("O8").Value = 20.58686954905052

$ws.Range("B9").Value = 15.93186519801911
$ws.Range("D9").Value = 4.985266487502726
$ws.Range("E9").Value = 21.30338721411833
$ws.Range("F9").Value = 22.82001584786747
$ws.Range("G9").Value = 26.20081779668133
$ws.Range("H9").Value = 13.6463952428357
$ws.Range("K9").Value = 12.50512246954439
$ws.Range("L9").Value = 7.993514397770542
$ws.Range("M9").Value = 14.64962813526714
$ws.Range("O9").Value = 20.43020994831769

$ws.Range("B10").Value = 16.10699867881871
$ws.Range("D10").Value = 5.113997448201102
$ws.Range("E10").Value = 21.16839460191105
$ws.Range("F10").Value = 22.80454370140825
$ws.Range("G10").Value = 26.16075747684527
$ws.Range("H10").Value = 13.58221873594255
$ws.Range("K10").Value = 12.88689558258375
$ws.Range("L10").Value = 8.047803243959056
$ws.Range("M10").Value = 14.70875626943316
$ws.Range("O10").Value = 20.33652327189481

$ws.Range("B11").Value = 16.18884957143231
$ws.Range("D11").Value = 5.170867463969938
$ws.Range("E11").Value = 21.11000700057677
$ws.Range("F11").Value = 22.80399368626024
$ws.Range("G11").Value = 26.15192836776273
$ws.Range("H11").Value = 13.55526399918228
$ws.Range("K11").Value = 13.05570423997449
$ws.Range("L11").Value = 8.073320687705804
$ws.Range("M11").Value = 14.73751377219675
$ws.Range("O11").Value = 20.2985728377899

$ws.Range("B12").Value = 16.22013667583196
$ws.Range("D12").Value = 5.192149309685409
$ws.Range("E12").Value = 21.08832956193629
$ws.Range("F12").Value = 22.80471773898774
$ws.Range("G12").Value = 26.14993909809724
$ws.Range("H12").Value = 13.5453788909319
$ws.Range("K12").Value = 13.11889074740554
$ws.Range("L12").Value = 8.083096549170548
$ws.Range("M12").Value = 14.74866537869355
$ws.Range("O12").Value = 20.28487487739512

$ws.Range("B13").Value = 16.21338580531795
$ws.Range("D13").Value = 5.1875773233818
$ws.Range("E13").Value = 21.09297897065187
$ws.Range("F13").Value = 22.804520354562
$ws.Range("G13").Value = 26.15030726043824
$ws.Range("H13").Value = 13.54749350139931
$ws.Range("K13").Value = 13.10531579666726
$ws.Range("L13").Value = 8.080986201715968
$ws.Range("M13").Value = 14.74625213170765
$ws.Range("O13").Value = 20.28779501708815

$ws.Range("B14").Value = 16.19141785567216
$ws.Range("D14").Value = 5.172623474815819
$ws.Range("E14").Value = 21.10821492318346
$ws.Range("F14").Value = 22.80403457581955
$ws.Range("G14").Value = 26.15173755687617
$ws.Range("H14").Value = 13.55444429052072
$ws.Range("K14").Value = 13.06091762750828
$ws.Range("L14").Value = 8.074122719702181
$ws.Range("M14").Value = 14.73842600882904
$ws.Range("O14").Value = 20.29743239966816

$ws.Range("B15").Value = 16.17799921801677
$ws.Range("D15").Value = 5.163430483024683
$ws.Range("E15").Value = 21.11760367944871
$ws.Range("F15").Value = 22.80385840613975
$ws.Range("G15").Value = 26.15279007490823
$ws.Range("H15").Value = 13.55874379294842
$ws.Range("K15").Value = 13.03362530107074
$ws.Range("L15").Value = 8.069933194875263
$ws.Range("M15").Value = 14.7336661992966
$ws.Range("O15").Value = 20.30342327628726

$ws.Range("B16").Value = 16.10169158153348
$ws.Range("D16").Value = 5.110245931518814
$ws.Range("E16").Value = 21.172271011942
$ws.Range("F16").Value = 22.80471016704036
$ws.Range("G16").Value = 26.1615238241166
$ws.Range("H16").Value = 13.58402535806081
$ws.Range("K16").Value = 12.87576235360893
$ws.Range("L16").Value = 8.046151708412538
$ws.Range("M16").Value = 14.70691384395218
$ws.Range("O16").Value = 20.33909751031385

$ws.Range("B17").Value = 16.05542277387483
$ws.Range("D17").Value = 5.077178316226993
$ws.Range("E17").Value = 21.20658019103195
$ws.Range("F17").Value = 22.80689426820119
$ws.Range("G17").Value = 26.16929053572462
$ws.Range("H17").Value = 13.60010838539509
$ws.Range("K17").Value = 12.77764507214855
$ws.Range("L17").Value = 8.031769260198834
$ws.Range("M17").Value = 14.69097452756372
$ws.Range("O17").Value = 20.36217940539531

$ws.Range("B18").Value = 16.02901687682793
$ws.Range("D18").Value = 5.0580001870346
$ws.Range("E18").Value = 21.22659839894529
$ws.Range("F18").Value = 22.80876130755763
$ws.Range("G18").Value = 26.17464188780172
$ws.Range("H18").Value = 13.60956970564277
$ws.Range("K18").Value = 12.72075512709017
$ws.Range("L18").Value = 8.023574465652526
$ws.Range("M18").Value = 14.68198206523329
$ws.Range("O18").Value = 20.37589478452812

$ws.Range("B19").Value = 16.02011245056742
$ws.Range("D19").Value = 5.051479907992706
$ws.Range("E19").Value = 21.23342513333387
$ws.Range("F19").Value = 22.80949836931376
$ws.Range("G19").Value = 26.17660549838929
$ws.Range("H19").Value = 13.61280934827607
$ws.Range("K19").Value = 12.70141615669676
$ws.Range("L19").Value = 8.020813332226599
$ws.Range("M19").Value = 14.67896766773805
$ws.Range("O19").Value = 20.38061396793713

$ws.Range("B20").Value = 16.06032693143132
$ws.Range("D20").Value = 5.080714903797227
$ws.Range("E20").Value = 21.20289849106933
$ws.Range("F20").Value = 22.80659855372128
$ws.Range("G20").Value = 26.1683722231479
$ws.Range("H20").Value = 13.59837450365339
$ws.Range("K20").Value = 12.78813725206584
$ws.Range("L20").Value = 8.033292300757285
$ws.Range("M20").Value = 14.69265318094312
$ws.Range("O20").Value = 20.3596768208781

$ws.Range("B21").Value = 16.19786262748644
$ws.Range("D21").Value = 5.177022746046831
$ws.Range("E21").Value = 21.10372802581695
$ws.Range("F21").Value = 22.80415196681936
$ws.Range("G21").Value = 26.15128067440123
$ws.Range("H21").Value = 13.55239393424459
$ws.Range("K21").Value = 13.07397876045534
$ws.Range("L21").Value = 8.076135666611883
$ws.Range("M21").Value = 14.74071767309437
$ws.Range("O21").Value = 20.29458338683677

$ws.Range("B22").Value = 16.28944144896785
$ws.Range("D22").Value = 5.238482881104263
$ws.Range("E22").Value = 21.04143557054717
$ws.Range("F22").Value = 22.80798646703725
$ws.Range("G22").Value = 26.14800380818235
$ws.Range("H22").Value = 13.52422023382298
$ws.Range("K22").Value = 13.25647695629974
$ws.Range("L22").Value = 8.104792220800922
$ws.Range("M22").Value = 14.77365355216932
$ws.Range("O22").Value = 20.25596484487053

$ws.Range("B23").Value = 16.24041666440138
$ws.Range("D23").Value = 5.20581943214009
$ws.Range("E23").Value = 21.07445211756682
$ws.Range("F23").Value = 22.80544318044178
$ws.Range("G23").Value = 26.14902975237768
$ws.Range("H23").Value = 13.53908528799172
$ws.Range("K23").Value = 13.15948119347225
$ws.Range("L23").Value = 8.089439389557027
$ws.Range("M23").Value = 14.75593764255478
$ws.Range("O23").Value = 20.27621669886467

$ws.Range("B24").Value = 16.05810915381559
$ws.Range("D24").Value = 5.079116532284731
$ws.Range("E24").Value = 21.20456207356598
$ws.Range("F24").Value = 22.80673034182983
$ws.Range("G24").Value = 26.16878463218693
$ws.Range("H24").Value = 13.59915772195961
$ws.Range("K24").Value = 12.78339523284235
$ws.Range("L24").Value = 8.032603503650297
$ws.Range("M24").Value = 14.69189372785909
$ws.Range("O24").Value = 20.36080685242359

$ws.Range("B25").Value = 15.86944620089139
$ws.Range("D25").Value = 4.936434295363527
$ws.Range("E25").Value = 21.35579045499871
$ws.Range("F25").Value = 22.83146651459988
$ws.Range("G25").Value = 26.22387143909537
$ws.Range("H25").Value = 13.67201440166022
$ws.Range("K25").Value = 12.36048984999863
$ws.Range("L25").Value = 7.974283045666369
$ws.Range("O25").Value = 20.4688377220426
